$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Atualizando o Método de Estimativa (coluna ESTIMATIVA INICIAL) para cada item do backlog
$ws.Range("E5").Value = "6 pontos`n"
$ws.Range("E6").Value = "3 pontos"
$ws.Range("E7").Value = "5 pontos"
$ws.Range("E8").Value = "6 pontos"
$ws.Range("E9").Value = "3 pontos"
$ws.Range("E10").Value = "5 pontos"
$ws.Range("E11").Value = "3 pontos"
$ws.Range("E12").Value = "4 pontos"

# Atualizando a seleção/visão ativa da planilha
$ws.Range("E11").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 4
